$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "slab_thick"
$ws.Range("B26").Value = 1.365
$ws.Range("C26").Value = "in"

$ws.Range("A26").Select()
